$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion text (cell A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$oldLine1 = [string][char]0x2705 + " 1000 Bs = 2.7 = 10002.97 pesos"
$newLine1 = [string][char]0x2705 + " 1000 Bs = 2.59 = 9590.57 pesos"

$oldLine2 = [string][char]0x2705 + " 10002.97 pesos = 2.69 = 942.8 Bs"
$newLine2 = [string][char]0x2705 + " 9590.57 pesos = 2.57 = 926.87 Bs"

$currentText = $wsHoja1.Range("A1").Value2
$newText = $currentText.Replace($oldLine1, $newLine1).Replace($oldLine2, $newLine2)
$wsHoja1.Range("A1").Value2 = $newText

# --- Sheet "tasas": update rate cells N10, O10, N12, O12 ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 385.9
$wsTasas.Range("O10").Value = 3701
$wsTasas.Range("N12").Value = 3725
$wsTasas.Range("O12").Value = 360
